$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User審查意見彙整")

$ws.Rows.Item(72).Delete()
$ws.Rows.Item(70).Delete()

$pc = $wb.PivotCaches().Item(1)
try {
    $pc.ChangeSource("User審查意見彙整!`$A`$1:`$K`$83")
    Write-Host ("SourceData after ChangeSource: " + $pc.SourceData)
} catch {
    Write-Host ("Error ChangeSource: " + $_)
}
$pc.Refresh()
Write-Host ("SourceData after refresh: " + $pc.SourceData)
